$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Estipular o tema e criar o Termo de Abertura do Projeto"
$ws.Range("B5").Value = "Configurar e Organizar o GitHub e Ferramenta de Gestão do Projeto`n"
$ws.Range("B20").Value = "Criar o banco de dados em seu modelo lógico e físico"
$ws.Range("B22").Value = "Fazer Bateria de Testes de Integração entre as 3 Camadas: Jogo, Login e Cadastro"

$ws.Range("B23").Select()
